$d = $word.ActiveDocument

# --- Edit 1: footer PAGE field cached result "2" -> "1" --------------------
# The default footer (index 1 = primary) has a PAGE field whose cached
# text result needs to change from "2" to "1".
$sec = $d.Sections.Item(1)
$ftr = $sec.Footers.Item(1)
$pageField = $ftr.Range.Fields.Item(1)
$pageResult = $pageField.Result
$pageResult.Find.Execute("2", $false, $false, $false, $false, $false, $true, 1, $false, "1", 2)

# --- Edit 2: header date line - prepend a new date entry -------------------
# "17.8.13(N) " -> "17.10.8(E) 17.8.13(N) " (insert a new date before the
# existing one). Track the insertion as a revision and accept it immediately
# so the new text lands in its own run (matching identical formatting to the
# original run) without Word merging it into neighbouring runs.
$hdr = $sec.Headers.Item(1)
$d.TrackRevisions = $true
$hdrRng = $hdr.Range
$found = $hdrRng.Find.Execute("17.8.13(N) ")
if ($found) {
    $hdrRng.InsertBefore("17.10.8(E) ")
}
$d.AcceptAllRevisions()
$d.TrackRevisions = $false
